# Insert a new record row at row 712 (pushing the existing rows 712:801
# down to 713:802) and populate it with the new price-report entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 712:801 down to 713:802, leaving a blank row 712 in place
# (format is inherited from the row above, matching column D's date style).
$ws.Rows(712).Insert()

# Populate the newly inserted row 712 with the new record.
$ws.Range("A712").Value = 3
$ws.Range("B712").Value = "Femacal de La Calera"
$ws.Range("C712").Value = "Coquimbo"
$ws.Range("D712").Value = 45124
$ws.Range("E712").Value = 5
$ws.Range("F712").Value = 100112032
$ws.Range("G712").Value = "Zapallo italiano"
$ws.Range("H712").Value = "Sin especificar"
$ws.Range("I712").Value = "Primera"
$ws.Range("J712").Value = 120
$ws.Range("K712").Value = 11000
$ws.Range("L712").Value = 11500
$ws.Range("M712").Value = 11271
$ws.Range("N712").Value = "`$/caja 60 unidades"
$ws.Range("O712").Value = "Región de Arica y Parinacota"
$ws.Range("P712").Value = 188
$ws.Range("Q712").Value = 60
$ws.Range("R712").Value = "Hortaliza"
